# v2.0 Review LOGIN and No comments
# Fill in row 7 of the "Version History" sheet with the new v2.0 entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Version History")

$ws.Range("A7").Value = "v2.0"
$ws.Range("B7").Value = "Hala Eldaly"
$ws.Range("C7").Value = "No Comments "
$ws.Range("D7").Value = Get-Date -Year 2025 -Month 5 -Day 13 -Hour 0 -Minute 0 -Second 0

$ws.Range("D8").Select()
